$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header H1 keeps its text ("Паспорт") but gets a new style: bold + centered + text format.
# Setting NumberFormat to Text ("@") before re-assigning the (unchanged) value upgrades the
# existing bold/centered style (fontId=1, alignment=center) to also include numFmtId 49,
# producing the new cellXfs entry used by H1.
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "Паспорт"

# Fill in the passport number for every data row (2-5), formatted as text so Excel stores it
# as a shared string rather than a number (keeps the literal "8673515150" exactly as typed).
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "8673515150"

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "8673515150"

$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "8673515150"

$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "8673515150"

# Move/save the active selection to H3, matching the workbook's last edited cell.
$ws.Range("H3").Select()
